$wb = $excel.ActiveWorkbook

function Move-CaseTypeColumns($ws, $maxRow) {
  # Rotates the header block B:E so that:
  #   old B (Group)       -> new C
  #   old C (Data Type)   -> new D
  #   old D (Description) -> new E
  #   old E (Label)       -> new B
  # Columns A (Case Property) and F (Deprecated) are untouched.
  for ($r = 1; $r -le $maxRow; $r++) {
    $b = $ws.Cells.Item($r,2).Value2  # Group
    $c = $ws.Cells.Item($r,3).Value2  # Data Type
    $d = $ws.Cells.Item($r,4).Value2  # Description
    $e = $ws.Cells.Item($r,5).Value2  # Label

    if ($e -eq "" -or $e -eq $null) { $ws.Cells.Item($r,2).Clear() } else { $ws.Cells.Item($r,2).Value = $e }
    if ($b -eq "" -or $b -eq $null) { $ws.Cells.Item($r,3).Clear() } else { $ws.Cells.Item($r,3).Value = $b }
    if ($c -eq "" -or $c -eq $null) { $ws.Cells.Item($r,4).Clear() } else { $ws.Cells.Item($r,4).Value = $c }
    if ($d -eq "" -or $d -eq $null) { $ws.Cells.Item($r,5).Clear() } else { $ws.Cells.Item($r,5).Value = $d }
  }
}

function Reset-ColumnWidthKeepCellStyle($ws, $col, $maxRow, $refStyleRange) {
  # ClearFormats() is the only reliable way to drop a stale per-column width
  # override, but it also blanks the style of every cell in the column - so
  # re-stamp each used cell with a normal data style afterwards.
  $refStyle = $refStyleRange.Style
  $ws.Columns.Item($col).ClearFormats()
  for ($r = 1; $r -le $maxRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    if ($cell.Value2 -ne "" -and $cell.Value2 -ne $null) {
      $cell.Style = $refStyle
    }
  }
}

$ws1 = $wb.Worksheets.Item(1)   # caseType1
$ws2 = $wb.Worksheets.Item(2)   # caseType1-vl (unaffected data-wise)
$ws3 = $wb.Worksheets.Item(3)   # caseType2

# --- caseType1: move the Label column from E to B (Group/Data Type/Description shift right) ---
Move-CaseTypeColumns $ws1 5

# The custom column width that used to belong to column D (Description, 10.72 chars)
# now belongs to column E, since Description moved there.
Reset-ColumnWidthKeepCellStyle $ws1 4 5 $ws1.Range("A1")
$ws1.Columns.Item(5).ColumnWidth = 9.86

# --- caseType2: same column rotation ---
Move-CaseTypeColumns $ws3 4

# caseType2 gains a custom width on column A (Case Property)
$ws3.Columns.Item(1).ColumnWidth = 11.86

# --- Active sheet / selection bookkeeping ---
# caseType1 becomes the active tab (was caseType2).
$ws1.Activate()
$ws1.Range("E19").Select()

$ws3.Range("F1").Select()

$ws2.Range("B6").Select()

# Re-activate caseType1 so it is left as the selected/active sheet.
$ws1.Activate()
